# Daily attendance processing - 2025-10-27 15:23:49
#
# Reorders the comma-separated list of editors in the "Recorded By" column (G)
# for every data row on the active sheet. The reprocessing job re-sorts each
# cell's editor list by (descending) name length, stably, and then reverses
# the result - this has the effect of moving short "System"/"system" style
# tokens toward the front of the list while leaving single-editor cells (and
# already-correctly-ordered cells) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row()
$lastRow = $firstRow + $usedRange.Rows.Count() - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    $rawParts = $text.Split(",")
    if ($rawParts.Length -le 1) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $sortedDesc = $parts | Sort-Object -Property Length -Descending

    $count = $sortedDesc.Length
    $reordered = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reordered += $sortedDesc[$i]
    }

    $newText = $reordered -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
